# Fix the LOQ4241 course-syllabus sheet: the A-column labels for rows 10-23
# had drifted out of sync with the B/C "content" columns (and a few pieces
# of content were duplicated/misplaced). Re-align everything to the correct
# label/value pairing, then drop the now-unused trailing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Objetivos: / Docentes responsáveis name (content swap)
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C10").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Rows.Item(10).RowHeight = 60

# Row 11: Objectives: (label only, unchanged)
$ws.Range("A11").Value = "Objectives:"
$ws.Rows.Item(11).RowHeight = 60

# Row 12: Docentes responsáveis: (label only, unchanged)
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13: Programa resumido: / Semestral
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: Short syllabus: (label only)
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Rows.Item(14).RowHeight = 60

# Row 15: Programa: / 01/01/2016
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2016"
$ws.Range("C15").Value = "01/01/2016"
$ws.Rows.Item(15).RowHeight = 120

# Row 16: Syllabus: (label only)
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Rows.Item(16).RowHeight = 120

# Row 17: Avaliação: (label only, default height now)
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).RowHeight = 15

# Row 18: Método: / Docentes responsáveis name (again)
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Rows.Item(18).RowHeight = 60

# Row 19: Critério: / Aulas expositivas ...
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."
$ws.Rows.Item(19).RowHeight = 60

# Row 20: Norma de recuperação: / Nota Final formula text
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A Nota Final do aluno será determinada segundo a seguinte equação: Nota Final = (Prova- Bimestral-1*0,4) + (Prova-Bimestral-2*0,4) + (Trabalho*0,2)"
$ws.Range("C20").Value = "A Nota Final do aluno será determinada segundo a seguinte equação: Nota Final = (Prova- Bimestral-1*0,4) + (Prova-Bimestral-2*0,4) + (Trabalho*0,2)"
$ws.Rows.Item(20).RowHeight = 60

# Row 21: Bibliografia: / Média aritmética ...
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Rows.Item(21).RowHeight = 120

# Row 22: Requisitos: (label only, content below moves out)
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).RowHeight = 15

# Row 23: requirement text only (no A label)
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

# Row 24 is now obsolete (its content moved up into row 23) - remove it.
$ws.Rows.Item(24).Delete()
